# Apply the "added new heart, ufo mask and watermelon tournament and icon for Ufo mask" change.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # tournaments
$ws2 = $wb.Worksheets.Item(2)   # quests

# Row 1 on both sheets has no cells at all (it is just an empty, formatted
# spacer row), which makes the engine shrink the sheet dimension to start at
# A2 as soon as the file is touched. Re-stamping A1 with the (already
# default) "Normal" style materialises the cell again without introducing a
# new cell style, keeping the dimension anchored at A1.
$ws1.Cells.Item(1, 1).Style = "Normal"
$ws2.Cells.Item(1, 1).Style = "Normal"

# ---------------------------------------------------------------------------
# 1. Sheet1 "tournaments": fill rows 156-164 with the new tournament entries.
# ---------------------------------------------------------------------------

# Stamp the standard row formatting (columns A:G) used by similar existing
# "kill" tournament rows (row 100 has style pattern A9/B13/C9/D9/E7/F9/G12)
# onto every new/refreshed row, then overwrite the values.
$ws1.Range("A100:G100").Copy()
$ws1.Range("A156:G161").PasteSpecial(-4122)
$ws1.Range("A162:G164").PasteSpecial(-4122)

# Rows 162-164 use a different icon style (s=11, as seen e.g. on G59/G65)
# instead of the s=12 used by rows 156-161.
$ws1.Range("G59").Copy()
$ws1.Range("G162:G164").PasteSpecial(-4122)

$newRows = @(
    @{ Row = 156; B = "TID_EVENT_TOURNAMENT_KILL_NORMAL_HEART";                          C = "kill";          D = 0; E = "valentines_heart";        G = "icon_heart" },
    @{ Row = 157; B = "TID_EVENT_TOURNAMENT_KILL_TIME_LIMIT_HEART";                      C = "kill";          D = 2; E = "valentines_heart";        G = "icon_heart" },
    @{ Row = 158; B = "TID_EVENT_TOURNAMENT_KILL_TIME_ATTACK_HEART";                     C = "kill";          D = 1; E = "valentines_heart";        G = "icon_heart" },
    @{ Row = 159; B = "TID_EVENT_TOURNAMENT_KILL_NORMAL_WATERMELON";                     C = "kill";          D = 0; E = "PreSummer_02_watermelon"; G = "icon_watermelon" },
    @{ Row = 160; B = "TID_EVENT_TOURNAMENT_KILL_TIME_LIMIT_WATERMELON";                 C = "kill";          D = 2; E = "PreSummer_02_watermelon"; G = "icon_watermelon" },
    @{ Row = 161; B = "TID_EVENT_TOURNAMENT_KILL_TIME_ATTACK_WATERMELON";                C = "kill";          D = 1; E = "PreSummer_02_watermelon"; G = "icon_watermelon" },
    @{ Row = 162; B = "TID_EVENT_TOURNAMENT_WEARING_KILL_TIME_LIMIT_ALIEN_DISGUISE";     C = "kill_equipped"; D = 2; E = "alien_mask";               G = "icon_ufo_disguise_humanoids" },
    @{ Row = 163; B = "TID_EVENT_TOURNAMENT_KILL_WEARING_NORMAL_ALIEN_DISGUISE";         C = "kill_equipped"; D = 0; E = "alien_mask";               G = "icon_ufo_disguise_humanoids" },
    @{ Row = 164; B = "TID_EVENT_TOURNAMENT_WEARING_KILL_TIME_ATTACK_ALIEN_DISGUISE";    C = "kill_equipped"; D = 1; E = "alien_mask";               G = "icon_ufo_disguise_humanoids" }
)

# Shared-string table order matters: set every brand-new "Definition id" (column
# B) first, in row order, then fill in the remaining columns (which all reuse
# strings that already exist elsewhere in the workbook, except for the final
# "alien_mask" value used by the three alien-disguise rows).
foreach ($r in $newRows) {
    $ws1.Cells.Item($r.Row, 1).Value = "<Definition>"
    $ws1.Cells.Item($r.Row, 2).Value = $r.B
}

foreach ($r in $newRows) {
    $ws1.Cells.Item($r.Row, 3).Value = $r.C
    $ws1.Cells.Item($r.Row, 4).Value = $r.D
    $ws1.Cells.Item($r.Row, 7).Value = $r.G
}

foreach ($r in $newRows) {
    $ws1.Cells.Item($r.Row, 5).Value = $r.E
}

# Row 162 has no data/format at all in column F (unlike the other new rows,
# which keep an empty, styled F cell).
$ws1.Range("F162").Clear()

# ---------------------------------------------------------------------------
# 2. Sheet1 column widths (E and F got narrower).
# ---------------------------------------------------------------------------
$ws1.Columns.Item(5).ColumnWidth = 33.736979166666664
$ws1.Columns.Item(6).ColumnWidth = 20.451822916666668

# ---------------------------------------------------------------------------
# 3. Update the remembered selections on both sheets (without permanently
#    changing which sheet is active - "tournaments" stays the active tab).
# ---------------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("D6:F6").Select()

$ws1.Activate()
$ws1.Range("B162").Select()
